$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D2 value to 0
$ws.Range("D2").Value = 0

# Delete rows 3 to 5 (shrinks used range down to A1:E2)
$ws.Range("A3:E5").Delete()
